$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.830.58"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "2.951.28"
$ws.Range("E3").Value = "  +2.86%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.40"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.61"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.564"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.633"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.60"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0899"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +5.91%  "
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.06"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").Value = "3.421.38"
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("D16").Value = "2.952.10"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.998"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "51.945.89"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.72"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.55"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +7.16%  "
$ws.Range("E21").Value = "  -2.19%  "
$ws.Range("D22").Value = "0.0₃0990"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.51"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "272.58"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.79"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  +11.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.39"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +3.27%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.45"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +19.18%  "
$ws.Range("E30").Value = "  +23.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.77"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.40"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +9.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "37.71"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -2.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "53.16"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0449"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -0.95%  "
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.39"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.82"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.86"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.65"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.74"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +5.32%  "
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.55"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("E46").Value = "  +2.15%  "
$ws.Range("D47").Value = "2.165.10"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.93"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -7.63%  "
$ws.Range("E49").Value = "  +3.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0337"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +5.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.930"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -3.07%  "
